$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..74 (value, count, date_min, date_max)
$data = @(
    ,@('?', 509, '000000', '17880828')
    ,@('Coimbra', 44, '15420000', '17420427')
    ,@('Lisboa', 44, '15460000', '17530612')
    ,@('Paris', 42, '15340815', '17590310')
    ,@('Roma', 41, '15400927', '17560709')
    ,@('Goa', 30, '15480000', '17360000')
    ,@('Évora', 24, '15660714', '17460405')
    ,@('Nancy', 15, '16270929', '17510827')
    ,@('Macau', 14, '16280000', '17490201')
    ,@('Avignon', 13, '16150926', '17370927')
    ,@('Landsberg', 13, '16230729', '17571009')
    ,@('Viena, Áustria', 13, '16271016', '17651018')
    ,@('Mechelen', 10, '0', '16720424')
    ,@('Nápoles', 10, '15560120', '17000311')
    ,@('Génova', 9, '16390915', '17680930')
    ,@('Pequim', 9, '17310308', '17730214')
    ,@('Tournai', 9, '15960705', '16960929')
    ,@('Bordeaux', 8, '16680920', '17431107')
    ,@('Alcalá', 6, '15550952', '15710624')
    ,@('Brno, República Checa', 6, '16060000', '17261009')
    ,@('Chieri', 5, '16280213', '16880120')
    ,@('Japão', 5, '15560000', '16260600')
    ,@('Cracóvia', 4, '16310816', '16760909')
    ,@('Lyon', 4, '16740927', '16930923')
    ,@('Palermo, Sicília', 4, '16181103', '16561022')
    ,@('Portugal', 4, '15480613', '16780000')
    ,@('Bolonha', 3, '16680000', '17250728')
    ,@('Hangchow', 3, '16270000', '16790621')
    ,@('Milão', 3, '16731021', '16901101')
    ,@('Toulouse', 3, '16800921', '17280930')
    ,@('Trier', 3, '16770821', '17271019')
    ,@('Arona', 2, '15960921', '16140824')
    ,@('Boémia', 2, '17290927', '17291010')
    ,@('Douai', 2, '15941109', '15990512')
    ,@('Ferrara', 2, '15550000', '15561200')
    ,@('Leoben', 2, '16871224', '16881012')
    ,@('Messina, Sicília', 2, '15820000', '17061207')
    ,@('Montmartre', 2, '15340815', '15340815')
    ,@('Saragoça', 2, '16740613', '16771115')
    ,@('Shiuchow', 2, '15910101', '16050000')
    ,@('Shiuchow, Chao-tcheou fou', 2, '15891100', '15891100')
    ,@('Vilnius', 2, '16180531', '16720811')
    ,@('Bahia', 1, '16410000', '16410000')
    ,@('China', 1, '17280105', '17280105')
    ,@('Colorno', 1, '17991116', '17991116')
    ,@('Courtrai', 1, '16440926', '16440926')
    ,@('Goa, Índia', 1, '15791200', '15791200')
    ,@('Japão (província)', 1, '17280523', '17280523')
    ,@('Krems', 1, '16641031', '16641031')
    ,@('Lima, Peru', 1, '15680711', '15680711')
    ,@('Lorette', 1, '15590425', '15590425')
    ,@('Mainz', 1, '16730717', '16730717')
    ,@('Manila', 1, '16721011', '16721011')
    ,@('Mazowsze (província)', 1, '17700813', '17700813')
    ,@('Milão (província)', 1, '16581105', '16581105')
    ,@('Nan-tch''ang', 1, '16080000', '16080000')
    ,@('Nanquim', 1, '16080300', '16080300')
    ,@('Novellara', 1, '16001101', '16001101')
    ,@('Novellario', 1, '16170121', '16170121')
    ,@('Ormuz', 1, '000000', '000000')
    ,@('Ozukio (noviciado)', 1, '15811100', '15811100')
    ,@('Polotsk', 1, '17860903', '17860903')
    ,@('Salamanca', 1, '16830419', '16830419')
    ,@('Shanghai', 1, '16100000', '16100000')
    ,@('Todos-os-Santos, Nagasaki', 1, '16070202', '16070202')
    ,@('Toulouse (província)', 1, '17500319', '17500319')
    ,@('Trenčín', 1, '17291027', '17291027')
    ,@('Valença', 1, '15610927', '15610927')
    ,@('Veneza', 1, '17180424', '17180424')
    ,@('Villaregio', 1, '15890406', '15890406')
    ,@('[Missão]', 1, '17310308', '17310308')
    ,@('província do Japão', 1, '16910113', '16910113')
    ,@('Índia', 1, '15530000', '15530000')
)

# Force text formatting on the date columns so numeric-looking strings
# like "000000" / "17880828" are stored as text, not numbers.
$ws.Range("C2:D75").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Original sheet had 77 data rows (through row 78); target only has 73
# (through row 74), so remove the trailing now-unused rows.
$ws.Range("A75:D78").EntireRow.Delete()

